$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.389.81"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.630.59"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'607.11"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'147.03"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("E10").Value = "  +6.05%  "
$ws.Range("E11").Value = "  -2.60%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "'27.25"
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "3.101.31"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "63.216.73"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "2.622.61"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "'11.60"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "'4.53"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("D20").Value = "'343.41"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").Value = "'5.73"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").Value = "'66.11"
$ws.Range("E24").Value = "  -2.14%  "
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("D26").Value = "'1.61"
$ws.Range("E26").Value = "  +4.46%  "
$ws.Range("D27").Value = "'9.05"
$ws.Range("E27").Value = "  +6.92%  "
$ws.Range("D28").Value = "'550.09"
$ws.Range("E28").Value = "  +2.95%  "
$ws.Range("E29").Value = "  -1.48%  "
$ws.Range("D30").Value = "'8.00"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0843"
$ws.Range("E33").Value = "  +4.28%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'1.76"
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("D35").Value = "'5.27"
$ws.Range("E35").Value = "  +3.50%  "
$ws.Range("D36").Value = "'168.84"
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'0.403"
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").Value = "'1.93"
$ws.Range("E39").Value = "  +5.87%  "
$ws.Range("D40").Value = "'18.97"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "'165.18"
$ws.Range("E42").Value = "  -5.01%  "
$ws.Range("D43").Value = "'39.91"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("D44").Value = "'3.77"
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("D45").Value = "'21.84"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").Value = "'0.0566"
$ws.Range("E46").Value = "  +0.64%  "
$ws.Range("D47").Value = "'0.627"
$ws.Range("E47").Value = "  -0.75%  "
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("D49").Value = "'1.96"
$ws.Range("E49").Value = "  +14.11%  "
$ws.Range("D50").Value = "'0.0953"
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").Value = "'18.67"
$ws.Range("E51").Value = "  -0.28%  "
